$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column widths (custom accuracy formatting pass) ---
# ColumnWidth input maps to stored OOXML width via (round(cw*6)+5)/6;
# use precise fractional inputs so the saved width lands on an exact integer.
$ws.Range("B:C").ColumnWidth = 7.166666666666667
$ws.Range("G:G").ColumnWidth = 7.166666666666667
$ws.Range("I:M").ColumnWidth = 7.166666666666667
$ws.Range("O:Q").ColumnWidth = 7.166666666666667
$ws.Range("T:T").ColumnWidth = 8.166666666666666
$ws.Range("X:X").ColumnWidth = 7.166666666666667
$ws.Range("AA:AD").ColumnWidth = 7.166666666666667
$ws.Range("AH:AH").ColumnWidth = 7.166666666666667

# --- Replace data rows 2-5 with new 1000-point dataset sample, drop old row 6 ---
$data = New-Object 'object[,]' 4,34
$data[0,0] = 45067.50694444445
$data[0,1] = 5.885
$data[0,2] = 1.942
$data[0,3] = 1.363
$data[0,4] = 7.767
$data[0,5] = 3.737
$data[0,6] = 0.96
$data[0,7] = 6.281
$data[0,8] = 2.168
$data[0,9] = 0.758
$data[0,10] = 1.218
$data[0,11] = 2.359
$data[0,12] = 5.493
$data[0,13] = 0.668
$data[0,14] = 0.884
$data[0,15] = 2.495
$data[0,16] = 1.484
$data[0,17] = 1.256
$data[0,18] = 0
$data[0,19] = 25.797
$data[0,20] = 5.514
$data[0,21] = 3.303
$data[0,22] = 3.938
$data[0,23] = 2.803
$data[0,24] = 0.246
$data[0,25] = 1.46
$data[0,26] = 1.634
$data[0,27] = 0.588
$data[0,28] = 2.048
$data[0,29] = 3.034
$data[0,30] = 0.186
$data[0,31] = 2.547
$data[0,32] = 0.594
$data[0,33] = 1.826
$data[1,0] = 45067.51388888889
$data[1,1] = 21.859
$data[1,2] = 15.593
$data[1,3] = 1.208
$data[1,4] = 45.641
$data[1,5] = 36.669
$data[1,6] = 15.605
$data[1,7] = 56.186
$data[1,8] = 24.464
$data[1,9] = 10.811
$data[1,10] = 16.225
$data[1,11] = 18.017
$data[1,12] = 20.235
$data[1,13] = 5.185
$data[1,14] = 15.399
$data[1,15] = 22.975
$data[1,16] = 13.271
$data[1,17] = 0.766
$data[1,18] = 0.468
$data[1,19] = 238.321
$data[1,20] = 44.985
$data[1,21] = 15.36
$data[1,22] = 30.725
$data[1,23] = 16.378
$data[1,24] = 2.11
$data[1,25] = 28.047
$data[1,26] = 13.145
$data[1,27] = 11.009
$data[1,28] = 13.731
$data[1,29] = 19.192
$data[1,30] = 0.173
$data[1,31] = 49.709
$data[1,32] = 8.204000000000001
$data[1,33] = 18.299
$data[2,0] = 45067.52083333334
$data[2,1] = 8.436
$data[2,2] = 5.87
$data[2,3] = 0.571
$data[2,4] = 17.18
$data[2,5] = 13.542
$data[2,6] = 5.614
$data[2,7] = 27.83
$data[2,8] = 8.971
$data[2,9] = 3.907
$data[2,10] = 5.863
$data[2,11] = 6.705
$data[2,12] = 7.828
$data[2,13] = 1.935
$data[2,14] = 5.474
$data[2,15] = 8.601000000000001
$data[2,16] = 4.922
$data[2,17] = 0.463
$data[2,18] = 0.102
$data[2,19] = 84.01000000000001
$data[2,20] = 16.983
$data[2,21] = 5.817
$data[2,22] = 11.7
$data[2,23] = 6.216
$data[2,24] = 0.785
$data[2,25] = 12.811
$data[2,26] = 4.89
$data[2,27] = 3.931
$data[2,28] = 5.167
$data[2,29] = 7.212
$data[2,30] = 0.141
$data[2,31] = 24.868
$data[2,32] = 2.958
$data[2,33] = 6.711
$data[3,0] = 45067.52777777778
$data[3,1] = 10.15
$data[3,2] = 7.31
$data[3,3] = 0.54
$data[3,4] = 21.24
$data[3,5] = 17.08
$data[3,6] = 7.23
$data[3,7] = 29.83
$data[3,8] = 11.4
$data[3,9] = 4.98
$data[3,10] = 7.51
$data[3,11] = 8.380000000000001
$data[3,12] = 9.41
$data[3,13] = 2.41
$data[3,14] = 7.1
$data[3,15] = 10.72
$data[3,16] = 6.19
$data[3,17] = 0.38
$data[3,18] = 0.18
$data[3,19] = 106.93
$data[3,20] = 20.99
$data[3,21] = 7.13
$data[3,22] = 14.35
$data[3,23] = 7.63
$data[3,24] = 0.98
$data[3,25] = 14.08
$data[3,26] = 6.12
$data[3,27] = 5.1
$data[3,28] = 6.42
$data[3,29] = 8.93
$data[3,30] = 0.12
$data[3,31] = 26.38
$data[3,32] = 3.8
$data[3,33] = 8.51

$ws.Range("A2:AH5").Value = $data

# Old row 6 no longer exists in the refreshed dataset
$ws.Rows("6").Delete()
